$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped/rotated rows ---
$row = @(13, 7100665, "Greece Super League 1", "Greece Super League 1", 45172.58333333334, "OFI Crete", "PAOK Salonika", 1, 0, "H", 4.75, 3.6, 1.75, 4.5, 3.8, 1.75, 0.75, 1.85, 2, 2.75, 2, 1.85, 3.5, -1, -1, 0.8500000000000001, -1, -1, 0.8500000000000001)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A15:AC15").Value = $arr

$row = @(15, 7100661, "Greece Super League 1", "Greece Super League 1", 45172.58333333334, "Aris Salonika", "Asteras Tripolis", 3, 2, "H", 1.8, 3.4, 4.75, 1.55, 3.8, 7, -1, 2, 1.85, 2.25, 1.825, 2.025, 0.55, -1, -1, 0, 0, 0.825, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A17:AC17").Value = $arr

$row = @(49, 6937192, "Greece Super League 1", "Greece Super League 1", 45221.41666666666, "Asteras Tripolis", "AEK Athens", 0, 3, "A", 6.5, 3.6, 1.6, 7.5, 4, 1.5, 1, 2.025, 1.825, 2.5, 1.925, 1.925, -1, -1, 0.5, -1, 0.825, 0.925, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A51:AC51").Value = $arr

$row = @(50, 6935736, "Greece Super League 1", "Greece Super League 1", 45221.41666666666, "Kifisias FC", "OFI Crete", 0, 0, "D", 3.3, 3.25, 2.25, 3.1, 3.4, 2.4, 0.25, 1.8, 2.05, 2.5, 1.825, 2.025, -1, 2.4, -1, 0.4, -0.5, -1, 1.025)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A52:AC52").Value = $arr

$row = @(110, 6935718, "Greece Super League 1", "Greece Super League 1", 45298.41666666666, "Panserraikos", "Kifisias FC", 1, 1, "D", 2.2, 3.3, 3.4, 1.8, 3.6, 4.5, -0.75, 2.05, 1.8, 2.5, 1.85, 2, -1, 2.6, -1, -1, 0.8, -1, 1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A112:AC112").Value = $arr

$row = @(111, 6937231, "Greece Super League 1", "Greece Super League 1", 45298.41666666666, "Volos NFC", "Atromitos Athinon", 1, 2, "A", 2.4, 3.3, 3, 2.25, 3.2, 3.5, -0.25, 1.875, 1.975, 2, 1.775, 2.1, -1, -1, 2.5, -1, 0.9750000000000001, 0.7749999999999999, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A113:AC113").Value = $arr

$row = @(122, 6936857, "Greece Super League 1", "Greece Super League 1", 45305.64583333334, "AEK Athens", "Panathinaikos", 2, 2, "D", 1.909, 3.5, 4.2, 2.15, 3.2, 3.5, -0.25, 1.85, 2, 2, 1.8, 2.05, -1, 2.2, -1, -0.5, 0.5, 0.8, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A124:AC124").Value = $arr

$row = @(123, 6937238, "Greece Super League 1", "Greece Super League 1", 45305.64583333334, "PAOK Salonika", "Giannina", 4, 0, "H", 1.111, 9, 23, 1.25, 6, 9, -1.75, 2.025, 1.825, 2.75, 1.8, 2.05, 0.25, -1, -1, 1.025, -1, 0.8, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A125:AC125").Value = $arr

$row = @(141, 6937247, "Greece Super League 1", "Greece Super League 1", 45326.52083333334, "AEK Athens", "Asteras Tripolis", 4, 2, "H", 1.285, 5.5, 12, 1.285, 5.75, 10, -1.5, 1.825, 2.025, 3, 2.025, 1.825, 0.2849999999999999, -1, -1, 0.825, -1, 1.025, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A143:AC143").Value = $arr

$row = @(142, 6937250, "Greece Super League 1", "Greece Super League 1", 45326.52083333334, "Giannina", "Lamia", 1, 4, "A", 2.3, 3.25, 3.25, 2.55, 2.875, 3.1, 0, 1.75, 2.125, 2, 1.85, 2, -1, -1, 2.1, -1, 1.125, 0.8500000000000001, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A144:AC144").Value = $arr

$row = @(168, 6937266, "Greece Super League 1", "Greece Super League 1", 45350.5, "Atromitos Athinon", "Lamia", 3, 1, "H", 2.3, 3.2, 3.1, 2.2, 3.3, 3.3, -0.25, 1.925, 1.925, 2.5, 2.025, 1.825, 1.2, -1, -1, 0.925, -1, 1.025, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A170:AC170").Value = $arr

$row = @(169, 6937268, "Greece Super League 1", "Greece Super League 1", 45350.5, "Panetolikos", "Olympiakos", 1, 2, "A", 8, 5, 1.363, 8.5, 5, 1.363, 1.25, 2.025, 1.825, 2.75, 1.85, 2, -1, -1, 0.363, 0.5125, -0.5, 0.425, -0.5)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A171:AC171").Value = $arr

$row = @(173, 6936863, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "OFI Crete", "Panathinaikos", 2, 2, "D", 8, 4.75, 1.4, 5.5, 4.75, 1.55, 1, 1.95, 1.9, 2.5, 1.85, 2, -1, 3.75, -1, 0.95, -1, 0.8500000000000001, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A175:AC175").Value = $arr

$row = @(174, 6935700, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "Panserraikos", "Asteras Tripolis", 2, 1, "H", 2.6, 3.2, 2.875, 2.25, 3.3, 3.3, -0.25, 1.925, 1.925, 2.25, 2, 1.85, 1.25, -1, -1, 0.925, -1, 1, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A176:AC176").Value = $arr

$row = @(175, 6937269, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "Aris Salonika", "AEK Athens", 3, 3, "D", 4.75, 3.75, 1.75, 6.5, 4.2, 1.5, 1, 2.05, 1.8, 2.5, 1.975, 1.875, -1, 3.2, -1, 1.05, -1, 0.9750000000000001, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A177:AC177").Value = $arr

$row = @(176, 6937271, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "Giannina", "Atromitos Athinon", 1, 1, "D", 2.45, 3.1, 3.1, 2, 3.3, 4, -0.5, 2.025, 1.825, 2.25, 1.85, 2, -1, 2.3, -1, -1, 0.825, -0.5, 0.5)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A178:AC178").Value = $arr

$row = @(177, 6937270, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "Olympiakos", "Volos NFC", 3, 0, "H", 1.125, 9, 19, 1.111, 9, 21, -2.25, 1.875, 1.975, 3.25, 2, 1.85, 0.111, -1, -1, 0.875, -1, -0.5, 0.425)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A179:AC179").Value = $arr

$row = @(179, 6937272, "Greece Super League 1", "Greece Super League 1", 45354.60416666666, "Lamia", "PAOK Salonika", 0, 2, "A", 7.5, 4.5, 1.444, 9.5, 5, 1.333, 1.5, 1.925, 1.925, 3, 1.95, 1.9, -1, -1, 0.333, -1, 0.925, -1, 0.8999999999999999)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A181:AC181").Value = $arr

# --- Add new rows 185-188 ---
$ws.Range("A184:AC184").Copy() | Out-Null
$ws.Range("A185:AC188").PasteSpecial(-4122) | Out-Null

$row = @(183, 7916134, "Greece Super League 1", "Greece Super League 1", 45367.5, "Asteras Tripolis", "Volos NFC", 0, 2, "A", 2, 3.3, 3.75, 2.55, 3.25, 2.875, 0, 1.8, 2.05, 2.25, 1.925, 1.925, -1, -1, 1.875, -1, 1.05, -0.5, 0.4625)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A185:AC185").Value = $arr

$row = @(184, 7916136, "Greece Super League 1", "Greece Super League 1", 45367.52083333334, "Kifisias FC", "OFI Crete", 0, 0, "D", 2.625, 3.25, 2.625, 2.2, 3.3, 3.4, -0.25, 1.875, 1.975, 2.5, 2, 1.85, -1, 2.3, -1, -0.5, 0.4875, -1, 0.8500000000000001)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A186:AC186").Value = $arr

$row = @(185, 7916135, "Greece Super League 1", "Greece Super League 1", 45367.60416666666, "Atromitos Athinon", "Giannina", 3, 2, "H", 2, 3.3, 3.75, 2.3, 3.3, 3.25, -0.25, 1.975, 1.875, 2.25, 1.975, 1.875, 1.3, -1, -1, 0.9750000000000001, -1, 0.9750000000000001, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A187:AC187").Value = $arr

$row = @(186, 7916137, "Greece Super League 1", "Greece Super League 1", 45367.625, "Panetolikos", "Panserraikos", 3, 0, "H", 2.05, 3.2, 3.75, 2.25, 3.2, 3.4, -0.25, 1.925, 1.925, 2.25, 2.05, 1.8, 1.25, -1, -1, 0.925, -1, 1.05, -1)
$arr = New-Object "object[,]" 1,29
for ($i=0; $i -lt 29; $i++) { $arr[0,$i] = $row[$i] }
$ws.Range("A188:AC188").Value = $arr
